# Auto-generated PowerShell Excel COM-interop script
# Applies updated attendance/ticket-count figures (column F) across the four
# sheets of the workbook, matching the commit 'Update gh-pages to output generated at 456a3b4'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
# F3: 3402 -> 3405
$ws.Cells.Item(3, 6).Value = 3405
# F8: 1894 -> 1896
$ws.Cells.Item(8, 6).Value = 1896
# F10: 1769 -> 1774
$ws.Cells.Item(10, 6).Value = 1774
# F11: 75 -> 76
$ws.Cells.Item(11, 6).Value = 76
# F16: 8592 -> 8597
$ws.Cells.Item(16, 6).Value = 8597
# F17: 214 -> 216
$ws.Cells.Item(17, 6).Value = 216
# F18: 1147 -> 1146
$ws.Cells.Item(18, 6).Value = 1146
# F19: 317 -> 318
$ws.Cells.Item(19, 6).Value = 318
# F24: 1203 -> 1204
$ws.Cells.Item(24, 6).Value = 1204
# F25: 1079 -> 1080
$ws.Cells.Item(25, 6).Value = 1080
# F26: 607 -> 608
$ws.Cells.Item(26, 6).Value = 608
# F27: 36 -> 37
$ws.Cells.Item(27, 6).Value = 37
# F31: 121 -> 122
$ws.Cells.Item(31, 6).Value = 122
# F33: 650 -> 651
$ws.Cells.Item(33, 6).Value = 651
# F36: 3617 -> 3618
$ws.Cells.Item(36, 6).Value = 3618
# F41: 528 -> 530
$ws.Cells.Item(41, 6).Value = 530
# F42: 135 -> 136
$ws.Cells.Item(42, 6).Value = 136
# F44: 723 -> 726
$ws.Cells.Item(44, 6).Value = 726

$ws = $wb.Worksheets.Item("演出")
# F4: 352 -> 355
$ws.Cells.Item(4, 6).Value = 355
# F10: 206 -> 207
$ws.Cells.Item(10, 6).Value = 207
# F11: 33 -> 35
$ws.Cells.Item(11, 6).Value = 35
# F12: 33 -> 35
$ws.Cells.Item(12, 6).Value = 35
# F20: 20 -> 21
$ws.Cells.Item(20, 6).Value = 21
# F21: 52 -> 54
$ws.Cells.Item(21, 6).Value = 54
# F23: 12 -> 13
$ws.Cells.Item(23, 6).Value = 13
# F24: 123 -> 124
$ws.Cells.Item(24, 6).Value = 124
# F25: 7001 -> 7006
$ws.Cells.Item(25, 6).Value = 7006
# F31: 11 -> 12
$ws.Cells.Item(31, 6).Value = 12
# F33: 58 -> 59
$ws.Cells.Item(33, 6).Value = 59
# F38: 36 -> 38
$ws.Cells.Item(38, 6).Value = 38

$ws = $wb.Worksheets.Item("本地生活")
# F4: 2111 -> 2113
$ws.Cells.Item(4, 6).Value = 2113
# F5: 1417 -> 1418
$ws.Cells.Item(5, 6).Value = 1418
# F6: 40 -> 41
$ws.Cells.Item(6, 6).Value = 41
# F8: 2258 -> 2259
$ws.Cells.Item(8, 6).Value = 2259
# F9: 9088 -> 9089
$ws.Cells.Item(9, 6).Value = 9089
# F10: 1367 -> 1369
$ws.Cells.Item(10, 6).Value = 1369
# F12: 29 -> 30
$ws.Cells.Item(12, 6).Value = 30

$ws = $wb.Worksheets.Item("全部类型")
# F3: 3402 -> 3405
$ws.Cells.Item(3, 6).Value = 3405
# F4: 2111 -> 2113
$ws.Cells.Item(4, 6).Value = 2113
# F5: 1417 -> 1418
$ws.Cells.Item(5, 6).Value = 1418
# F6: 2258 -> 2259
$ws.Cells.Item(6, 6).Value = 2259
# F7: 1367 -> 1369
$ws.Cells.Item(7, 6).Value = 1369
# F9: 29 -> 30
$ws.Cells.Item(9, 6).Value = 30
# F16: 8592 -> 8597
$ws.Cells.Item(16, 6).Value = 8597
# F17: 214 -> 216
$ws.Cells.Item(17, 6).Value = 216
# F18: 317 -> 318
$ws.Cells.Item(18, 6).Value = 318
# F22: 607 -> 608
$ws.Cells.Item(22, 6).Value = 608
# F23: 36 -> 37
$ws.Cells.Item(23, 6).Value = 37
# F26: 206 -> 207
$ws.Cells.Item(26, 6).Value = 207
# F27: 33 -> 35
$ws.Cells.Item(27, 6).Value = 35
# F30: 650 -> 651
$ws.Cells.Item(30, 6).Value = 651
# F35: 3617 -> 3618
$ws.Cells.Item(35, 6).Value = 3618
# F38: 528 -> 530
$ws.Cells.Item(38, 6).Value = 530
# F39: 52 -> 54
$ws.Cells.Item(39, 6).Value = 54
# F41: 723 -> 726
$ws.Cells.Item(41, 6).Value = 726
